# Insert a new weekly price record for "Femacal de La Calera - Cebollín"
# as row 520, pushing the existing rows 520:560 down to 521:561.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 520 (shifts 520:560 -> 521:561).
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(520, 1).Value = 3
$ws.Cells.Item(520, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(520, 3).Value = "Coquimbo"
$ws.Cells.Item(520, 4).Value = 44826
$ws.Cells.Item(520, 5).Value = 5
$ws.Cells.Item(520, 6).Value = 100112037
$ws.Cells.Item(520, 7).Value = "Cebollín"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 250
$ws.Cells.Item(520, 11).Value = 4500
$ws.Cells.Item(520, 12).Value = 5000
$ws.Cells.Item(520, 13).Value = 4740
$ws.Cells.Item(520, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(520, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(520, 16).Value = 132
$ws.Cells.Item(520, 17).Value = 36
$ws.Cells.Item(520, 18).Value = "Hortaliza"
